$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.29025
$ws.Range("H2").Value = 6.870749999999999
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.3777296666666667
$ws.Range("N2").Value = 1.133189
$ws.Range("O2").Value = 0.2121489992374768
$ws.Range("P2").Value = 0.2121489992374768
$ws.Range("Q2").Value = 0.8650953690833333
$ws.Range("R2").Value = 7.785858321749999
$ws.Range("S2").Value = 0.2121489992374768
$ws.Range("T2").Value = 0.2121489992374768

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.29025
$ws.Range("H3").Value = 6.870749999999999
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.092289666666667
$ws.Range("N3").Value = 3.276869
$ws.Range("O3").Value = 0.6134761976883921
$ws.Range("P3").Value = 0.6134761976883921
$ws.Range("Q3").Value = 2.501616409083333
$ws.Range("R3").Value = 22.51454768175
$ws.Range("S3").Value = 0.6134761976883921
$ws.Range("T3").Value = 0.6134761976883921

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.29025
$ws.Range("H4").Value = 6.870749999999999
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.310473
$ws.Range("N4").Value = 0.931419
$ws.Range("O4").Value = 0.174374803074131
$ws.Range("P4").Value = 0.174374803074131
$ws.Range("Q4").Value = 0.71106078825
$ws.Range("R4").Value = 6.399547094249999
$ws.Range("S4").Value = 0.174374803074131
$ws.Range("T4").Value = 0.174374803074131
